# Jimenez 2019 pattern generator: re-generated the AATT-repeat example
# strings (col C) and their generator labels / coincidence lists (cols B, E)
# for rows 2-48 of Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "o1980, o1988, o1992b, o2008"
$ws.Range("C3").Value = "ATATATTAAT"
$ws.Range("C4").Value = "ATATATAAAT"
$ws.Range("E4").Value = "o1980, o1988, o1992b"
$ws.Range("C5").Value = "ATATTAATAT"
$ws.Range("E5").Value = "o1980, o1988, o2008"
$ws.Range("C6").Value = "TAATATATAT"
$ws.Range("E6").Value = "o1980, o1988, o1992b, o2008"
$ws.Range("C7").Value = "ATAAATATAT"
$ws.Range("E7").Value = "o1980, o1988, o1992b"
$ws.Range("E8").Value = "o1980, o1988, o1992b"
$ws.Range("C9").Value = "ATATAAATAT"
$ws.Range("E9").Value = "o1980, o1988, o1992b"
$ws.Range("C10").Value = "ATTAATATAT"
$ws.Range("E10").Value = "o1980, o1988, o2008"
$ws.Range("C11").Value = "ATAAATAAAT"
$ws.Range("E11").Value = "o1980, o1988, o1992b"
$ws.Range("C12").Value = "AAAAATATAT"
$ws.Range("E12").Value = "o1980, o1988, o1992b"
$ws.Range("C13").Value = "AAATATAAAT"
$ws.Range("E13").Value = "o1980, o1988, o1992b"
$ws.Range("B14").Value = "o1980"
$ws.Range("C14").Value = "AAATTAATAT"
$ws.Range("E14").Value = "o1980, o1988"
$ws.Range("B15").Value = "o1980"
$ws.Range("C15").Value = "TAATAAATAT"
$ws.Range("E15").Value = "o1980, o1992b"
$ws.Range("B16").Value = "o1980"
$ws.Range("C16").Value = "AATAATATAT"
$ws.Range("E16").Value = "o1980, o1988, o1992b, o2008"
$ws.Range("B17").Value = "o1980"
$ws.Range("C17").Value = "ATTAAAATAT"
$ws.Range("E17").Value = "o1980"
$ws.Range("E18").Value = "o1980, o1988, o1992b"
$ws.Range("B19").Value = "o1980"
$ws.Range("C19").Value = "ATATAATAAT"
$ws.Range("E19").Value = "o1980, o1988, o1992b, o2008"
$ws.Range("C20").Value = "TAATATAAAT"
$ws.Range("E20").Value = "o1980, o1988, o1992b"
$ws.Range("C21").Value = "ATAATAATAT"
$ws.Range("E21").Value = "o1980, o1988"
$ws.Range("C22").Value = "AAATATTAAT"
$ws.Range("E22").Value = "o1980, o1988"
$ws.Range("C23").Value = "ATATAAAAAT"
$ws.Range("E23").Value = "o1980, o1988, o1992b"
$ws.Range("C24").Value = "ATAAATTAAT"
$ws.Range("E24").Value = "o1980"
$ws.Range("C25").Value = "ATAAAAATAT"
$ws.Range("E25").Value = "o1980, o1988"
$ws.Range("C26").Value = "AAATAAATAT"
$ws.Range("E26").Value = "o1980, o1988, o1992b"
$ws.Range("B27").Value = "o2008"
$ws.Range("C27").Value = "TAATTAATAT"
$ws.Range("E27").Value = "o2008"
$ws.Range("B28").Value = "o2008"
$ws.Range("C28").Value = "TAATATTAAT"
$ws.Range("E28").Value = "o2008"
$ws.Range("B29").Value = "o1988"
$ws.Range("C29").Value = "ATTAATAAAT"
$ws.Range("E29").Value = "o1988"
$ws.Range("B30").Value = "o2008"
$ws.Range("C30").Value = "ATTAATTAAT"
$ws.Range("E30").Value = "o2008"
$ws.Range("B31").Value = "o2008"
$ws.Range("C31").Value = "ATATTATAAT"
$ws.Range("E31").Value = "o2008"
$ws.Range("C32").Value = "ATAAAATAAT"
$ws.Range("B33").Value = "o1980"
$ws.Range("C33").Value = "AAATAAAAAT"
$ws.Range("E33").Value = "o1980, o1988, o1992b"
$ws.Range("B34").Value = "o1980"
$ws.Range("C34").Value = "AAATAATAAT"
$ws.Range("E34").Value = "o1980, o1988, o1992b"
$ws.Range("B35").Value = "o1980"
$ws.Range("C35").Value = "AATAAAATAT"
$ws.Range("E35").Value = "o1980"
$ws.Range("C36").Value = "AAAAATAAAT"
$ws.Range("E36").Value = "o1980, o1988, o1992b"
$ws.Range("C37").Value = "AAAAATTAAT"
$ws.Range("E37").Value = "o1980"
$ws.Range("C38").Value = "ATAAAAAAAT"
$ws.Range("E38").Value = "o1980, o1988"
$ws.Range("C39").Value = "TAAAAAATAT"
$ws.Range("E39").Value = "o1980"
$ws.Range("C40").Value = "AAAATAATAT"
$ws.Range("E40").Value = "o1980, o1988"
$ws.Range("C41").Value = "AAAAAAATAT"
$ws.Range("E41").Value = "o1980, o1988"
$ws.Range("C42").Value = "TAATAAAAAT"
$ws.Range("E42").Value = "o1980, o1992b"
$ws.Range("B43").Value = "o1988"
$ws.Range("C43").Value = "AATAATAAAT"
$ws.Range("E43").Value = "o1988, o1992b"
$ws.Range("C44").Value = "TAAAATAAAT"
$ws.Range("E44").Value = "o1980, o1988, o1992b"
$ws.Range("B45").Value = "o2008"
$ws.Range("C45").Value = "TAATAATAAT"
$ws.Range("E45").Value = "o2008"
$ws.Range("C46").Value = "TATAATTAAT"
$ws.Range("B47").Value = "o2008"
$ws.Range("C47").Value = "TAATTATAAT"
$ws.Range("E47").Value = "o2008"
$ws.Range("C48").Value = "AATAATTAAT"
